$wb = $excel.ActiveWorkbook

# --- GET_ClientCodeAssays: update test client code value, move selection ---
$wsClientCodeAssays = $wb.Worksheets.Item("GET_ClientCodeAssays")
$wsClientCodeAssays.Range("B3").Value = "CDF"

# --- GET_VerifiedSamples: update prefix/sampleID test values (two test blocks) ---
$wsVerifiedSamples = $wb.Worksheets.Item("GET_VerifiedSamples")
$wsVerifiedSamples.Range("B3").Value = "CDF"
$wsVerifiedSamples.Range("C3").Value = "CDF33"
$wsVerifiedSamples.Range("B8").Value = "CDF"
$wsVerifiedSamples.Range("C8").Value = "CDF33"

# --- Update each sheet's remembered selection (cursor position) ---
$wsAutoApproval = $wb.Worksheets.Item("GET_AutoApproval")
$wsAutoApproval.Range("F31").Select()

$wsClientCodeAssays.Range("B19").Select()

$wsVerifiedSamples.Range("B3").Select()

$wsCreateOrder = $wb.Worksheets.Item("POST_CreateOrder")
$wsCreateOrder.Range("E30").Select()

# --- Make GET_AutoApproval the active (front-most) sheet/tab ---
$wsAutoApproval.Activate()
